$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.841.90"
$ws.Range("E2").Value = "  -3.12%  "
$ws.Range("D3").Value = "3.853.43"
$ws.Range("E3").Value = "  -3.37%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.997"
$ws.Range("E4").Value = "  -0.26%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "597.93"
$ws.Range("E5").Value = "  -0.08%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "169.39"
$ws.Range("E6").Value = "  +5.29%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.666"
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("E8").Value = "  +0.17%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.742"
$ws.Range("E9").Value = "  -0.95%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.176"
$ws.Range("E10").Value = "  +4.71%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "53.08"
$ws.Range("E11").Value = "  -1.47%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000321"
$ws.Range("E12").Value = "  +0.87%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "11.25"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "4.462.38"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "21.17"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.845.41"
$ws.Range("E16").Value = "  -3.57%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "13.85"
$ws.Range("E17").Value = "  -1.39%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.20"
$ws.Range("E18").Value = "  -5.82%  "
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "70.519.85"
$ws.Range("E20").Value = "  -3.11%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "436.39"
$ws.Range("E21").Value = "  +0.21%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.76"
$ws.Range("E22").Value = "  -1.17%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "94.18"
$ws.Range("E23").Value = "  -1.76%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.28"
$ws.Range("E24").Value = "  -4.07%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "13.81"
$ws.Range("E25").Value = "  -3.08%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "11.56"
$ws.Range("E26").Value = "  +3.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "3.99"
$ws.Range("E27").Value = "  -7.96%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "5.97"
$ws.Range("E28").Value = "  +0.21%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "10.52"
$ws.Range("E29").Value = "  +1.29%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "34.92"
$ws.Range("E30").Value = "  -3.77%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.18"
$ws.Range("E31").Value = "  +4.92%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "13.46"
$ws.Range("E32").Value = "  -1.86%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "47.79"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  -3.88%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "68.81"
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("D36").Value = "0.0₃0977"
$ws.Range("E36").Value = "  +8.09%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "630.60"
$ws.Range("E37").Value = "  -5.44%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.435"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").Value = "  +0.25%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.996"
$ws.Range("E41").Value = "  -0.58%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.23"
$ws.Range("E42").Value = "  -3.00%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.91"
$ws.Range("E43").Value = "  +10.82%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.21"
$ws.Range("E44").Value = "  +21.67%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0469"
$ws.Range("E45").Value = "  -4.02%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "9.95"
$ws.Range("E46").Value = "  -5.94%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.143"
$ws.Range("E47").Value = "  -4.09%  "
$ws.Range("D48").Value = "2.891.76"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.80"
$ws.Range("E49").Value = "  -16.27%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "3.29"
$ws.Range("E50").Value = "  -4.72%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.000278"
$ws.Range("E51").Value = "  +3.26%  "
